# Apply edits to MyTestdata.xlsx as described by the target diff:
#  - rename sheet "loginpage" -> "LOGIN"
#  - update test data rows (TestCaseID shift, Run column, new Name "xyz",
#    new Data_1/Data_2 columns with email + password, email hyperlinks)
#  - widen column B, add a width for new column D
#  - update the active selection on the LOGIN sheet and on Sheet2

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("loginpage")
$ws.Name = "LOGIN"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 21.140625
$ws.Columns.Item(4).ColumnWidth = 24.140625

# --- Row 2 ---
$ws.Range("A2").Value = "anil"
$ws.Range("B2").Value = "'TestCase_001"
$ws.Range("C2").Value = "Y"

# --- Row 3 ---
$ws.Range("A3").Value = "anil"
$ws.Range("B3").Value = "TestCase_002"
$ws.Range("C3").Value = "No"

# --- Row 4 ---
$ws.Range("A4").Value = "anil"
$ws.Range("B4").Value = "TestCase_003"
$ws.Range("C4").Value = "No"

# --- Row 5 ---
$ws.Range("A5").Value = "anil"
$ws.Range("B5").Value = "TestCase_004"
$ws.Range("C5").Value = "No"

# --- Row 6 ---
$ws.Range("A6").Value = "xyz"
$ws.Range("B6").Value = "TestCase_005"
$ws.Range("C6").Value = "No"

# --- Row 7 ---
$ws.Range("A7").Value = "xyz"
$ws.Range("B7").Value = "TestCase_006"
$ws.Range("C7").Value = "No"

# --- Row 8 ---
$ws.Range("A8").Value = "xyz"
$ws.Range("B8").Value = "TestCase_007"
$ws.Range("C8").Value = "No"

# --- Row 9 ---
$ws.Range("A9").Value = "xyz"
$ws.Range("B9").Value = "TestCase_008"
$ws.Range("C9").Value = "No"

# --- Row 10 ---
$ws.Range("A10").Value = "xyz"
$ws.Range("B10").Value = "TestCase_009"
$ws.Range("C10").Value = "No"

# --- Row 11 ---
$ws.Range("A11").Value = "xyz"
$ws.Range("B11").Value = "TestCase_010"
$ws.Range("C11").Value = "No"

# --- New Data_1 / Data_2 columns (D/E) with email hyperlinks ---
$ws.Range("E2").Value = "password"
$ws.Range("E3").Value = "password"

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:prathap@steedserv.com")
$ws.Range("D2").Value = "prathap@steedserv.com"
$ws.Range("D2").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:prathap2@steedserv.com")
$ws.Range("D3").Value = "prathap2@steedserv.com"

# --- Selections ---
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Range("F6").Select()

[void]$ws.Activate()
[void]$ws.Range("E4").Select()
